$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in column C
$ws.Range("C2").Value = 48
$ws.Range("C3").Value = 37
$ws.Range("C4").Value = 15

# Update the view: scroll back to top-left A1 and change selection to E9
$ws.Range("E9").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
